$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$unionRange = $null
for ($r = 2; $r -le 16; $r += 2) {
    $r2 = $r + 1
    $ws.Range("I$r").Formula = "=ROUND(AVERAGE(H$r,H$r2),4)"
    $pairRange = $ws.Range("I$r`:I$r2")
    $pairRange.Merge() | Out-Null
    if ($unionRange -eq $null) {
        $unionRange = $pairRange
    } else {
        $unionRange = $ws.Application.Union($unionRange, $pairRange)
    }
}
$unionRange.HorizontalAlignment = -4108
$unionRange.VerticalAlignment = -4108
